$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = '2026-02-05 10:07:06'
$ws.Range("E3").Value = '2026-02-05 10:07:08'
$ws.Range("H3").NumberFormat = "@"
$ws.Range("H3").Value = '67%'
$ws.Range("K3").Value = '0.4 MJ/m2'
$ws.Range("L3").Value = '35.6 km/h - 124º 8:58 TU'
$ws.Range("O3").Value = '-2.1 °C'
$ws.Range("E4").Value = '2026-02-05 10:07:11'
$ws.Range("E5").Value = '2026-02-05 10:07:13'
$ws.Range("E6").Value = '2026-02-05 10:07:15'
$ws.Range("E7").Value = '2026-02-05 10:07:18'
$ws.Range("I7").Value = '0.2 mm'
$ws.Range("J7").Value = '994.2 hPa'
$ws.Range("K7").Value = '0.2 MJ/m2'
$ws.Range("M7").Value = '11.0 °C 8:50 TU'
$ws.Range("E8").Value = '2026-02-05 10:07:20'
$ws.Range("E9").Value = '2026-02-05 10:07:23'
$ws.Range("E10").Value = '2026-02-05 10:07:25'
$ws.Range("E11").Value = '2026-02-05 10:07:27'
$ws.Range("E12").Value = '2026-02-05 10:07:30'
$ws.Range("H12").NumberFormat = "@"
$ws.Range("H12").Value = '93%'
$ws.Range("I12").Value = '0.8 mm'
$ws.Range("K12").Value = '0.2 MJ/m2'
$ws.Range("M12").Value = '8.5 °C 9:29 TU'
$ws.Range("O12").Value = '7.1 °C'
$ws.Range("E13").Value = '2026-02-05 10:07:32'
$ws.Range("E14").Value = '2026-02-05 10:07:34'
$ws.Range("H14").NumberFormat = "@"
$ws.Range("H14").Value = '66%'
$ws.Range("I14").Value = '0.1 mm'
$ws.Range("M14").Value = '-2.1 °C 1:31 TU'
$ws.Range("N14").Value = '-3.3 °C 0:38 TU'
$ws.Range("O14").Value = '-2.7 °C'
$ws.Range("E15").Value = '2026-02-05 10:07:37'
$ws.Range("E16").Value = '2026-02-05 10:07:39'
$ws.Range("E17").Value = '2026-02-05 10:07:42'
$ws.Range("J17").Value = '997.8 hPa'
$ws.Range("K17").Value = '0.3 MJ/m2'
$ws.Range("M17").Value = '2.1 °C 8:59 TU'
$ws.Range("O17").Value = '0.1 °C'
$ws.Range("E18").Value = '2026-02-05 10:07:44'
$ws.Range("E19").Value = '2026-02-05 10:07:47'
$ws.Range("E20").Value = '2026-02-05 10:07:49'
$ws.Range("E21").Value = '2026-02-05 10:07:52'
$ws.Range("E22").Value = '2026-02-05 10:07:54'
$ws.Range("E23").Value = '2026-02-05 10:07:57'
$ws.Range("E24").Value = '2026-02-05 10:07:59'
$ws.Range("E25").Value = '2026-02-05 10:08:02'
$ws.Range("E26").Value = '2026-02-05 10:08:05'
$ws.Range("E27").Value = '2026-02-05 10:08:07'
$ws.Range("J27").Value = '993.6 hPa'
$ws.Range("K27").Value = '0.5 MJ/m2'
$ws.Range("M27").Value = '6.5 °C 8:59 TU'
$ws.Range("O27").Value = '4.3 °C'
$ws.Range("E28").Value = '2026-02-05 10:08:10'
$ws.Range("E29").Value = '2026-02-05 10:08:12'
$ws.Range("E30").Value = '2026-02-05 10:08:15'
$ws.Range("E31").Value = '2026-02-05 10:08:17'
$ws.Range("G31").Value = '1 cm'
$ws.Range("H31").NumberFormat = "@"
$ws.Range("H31").Value = '95%'
$ws.Range("J31").Value = '996.8 hPa'
$ws.Range("M31").Value = '4.9 °C 8:55 TU'
$ws.Range("O31").Value = '3.5 °C'
$ws.Range("E32").Value = '2026-02-05 10:08:20'
$ws.Range("E33").Value = '2026-02-05 10:08:22'
$ws.Range("E34").Value = '2026-02-05 10:08:24'
$ws.Range("E35").Value = '2026-02-05 10:08:27'
$ws.Range("K35").Value = '0.2 MJ/m2'
$ws.Range("M35").Value = '-3.0 °C 8:26 TU'
$ws.Range("O35").Value = '-3.7 °C'
$ws.Range("E36").Value = '2026-02-05 10:08:29'
$ws.Range("I36").Value = '1.4 mm'
$ws.Range("J36").Value = '995.4 hPa'
$ws.Range("K36").Value = '0.1 MJ/m2'
$ws.Range("M36").Value = '6.7 °C 8:28 TU'
$ws.Range("O36").Value = '5.5 °C'
